$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab/title (and workbook.xml sheet name) to reflect new "through" date
$ws.Name = "Through 2021-12-02"

# Update the "December" row label to reflect new "through" date
$ws.Range("A13").Value = "December (through 12-02)"

# Update the December (row 13) figures
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 8
$ws.Range("H13").Value = 14

# Update the Total (row 14) figures
$ws.Range("C14").Value = 569
$ws.Range("D14").Value = 827
$ws.Range("E14").Value = 689
$ws.Range("F14").Value = 536
$ws.Range("G14").Value = 1272
$ws.Range("H14").Value = 1658
